$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename Sheet1 -> Register, add a new sheet "TryEditor" right after it
# ---------------------------------------------------------------------------
$register = $wb.Worksheets.Item(1)
$register.Name = "Register"

$tryEditor = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $register)
$tryEditor.Name = "TryEditor"

# ---------------------------------------------------------------------------
# 2. Update the wording of the scenario names on the Register sheet
# ---------------------------------------------------------------------------
$register.Range("A2").Value = "validate the error message displayed when the username field is left empty"
$register.Range("A3").Value = "validate the error message displayed when the password field is left empty"
$register.Range("A4").Value = " validate the error message displayed when the confirm password field is left empty"
$register.Range("A5").Value = "validate the error message displayed when all the fields  left empty`t"

# A2's format is tweaked directly (its border is removed and the font size is
# pinned to 10) - this is what makes it diverge from the other scenario cells.
$register.Range("A2").Borders.LineStyle = 0
$register.Range("A2").Font.Size = 10

# ---------------------------------------------------------------------------
# 3. Populate the new TryEditor sheet
# ---------------------------------------------------------------------------
$tryEditor.Range("A1").Value = "scenario"
$tryEditor.Range("B1").Value = "Input"
$tryEditor.Range("C1").Value = "output"

$tryEditor.Range("A2").Value = "Verify that user is able to see output for valid python code"
$tryEditor.Range("B2").Value = 'print("Hello! Welcome")'
$tryEditor.Range("C2").Value = "Hello! Welcome"

$tryEditor.Range("A3").Value = "Verify that user receives error for invalid python code"
$tryEditor.Range("B3").Value = "hi"
$tryEditor.Range("C3").Value = "NameError: name 'hi' is not defined on line 1"

$tryEditor.Range("A4").Value = "Verify that user receives error when click on Run button without entering code"

# Column widths (closest achievable match to 35.5 / 23.5 / 36.63 characters
# given the engine's column-width quantisation)
$tryEditor.Columns.Item(1).ColumnWidth = 34.6667
$tryEditor.Columns.Item(2).ColumnWidth = 22.6667
$tryEditor.Columns.Item(3).ColumnWidth = 35.8

# ---------------------------------------------------------------------------
# 4. Formatting - build each required look on a helper cell once, then stamp
#    it onto the remaining cells with Copy + PasteSpecial(Formats) so the
#    style/font tables don't balloon with duplicate entries.
# ---------------------------------------------------------------------------

# --- "plain" look: thin border, no fill, Arial (theme colour), bottom
#     vertical align - used for header row + a handful of other cells.
$plain = $tryEditor.Range("A1")
$plain.Borders.LineStyle = 1
$plain.Font.ThemeColor = 1
$plain.Font.Name = "Arial"
$plain.VerticalAlignment = -4107
$plain.Copy()

$tryEditor.Range("B1").PasteSpecial(-4122)
$tryEditor.Range("C1").PasteSpecial(-4122)
$tryEditor.Range("B2").PasteSpecial(-4122)
$tryEditor.Range("B3").PasteSpecial(-4122)
$tryEditor.Range("C3").PasteSpecial(-4122)
$tryEditor.Range("C4").PasteSpecial(-4122)

# --- "code input" look: thin border, white fill, Calibri 12, top align,
#     wrap text - used for the scenario / code-input column.
$codeFont = $tryEditor.Range("A2")
$codeFont.Borders.LineStyle = 1
$codeFont.Interior.ColorIndex = 2
$codeFont.Font.Name = "Calibri"
$codeFont.Font.ThemeColor = 1
$codeFont.Font.Size = 12
$codeFont.WrapText = $true
$codeFont.VerticalAlignment = -4160
$codeFont.ShrinkToFit = $false
$codeFont.Copy()

$tryEditor.Range("A3").PasteSpecial(-4122)
$tryEditor.Range("A4").PasteSpecial(-4122)

# --- "output" look: thin border, white fill, Arial 12, bottom align
$outputFont = $tryEditor.Range("C2")
$outputFont.Borders.LineStyle = 1
$outputFont.Interior.ColorIndex = 2
$outputFont.Font.Name = "Arial"
$outputFont.Font.ThemeColor = 1
$outputFont.Font.Size = 12
$outputFont.VerticalAlignment = -4107

# --- "input value, top aligned" look: thin border, white fill, Arial (theme
#     colour), top align - used for the second code snippet input cell.
$inputTop = $tryEditor.Range("B4")
$inputTop.Borders.LineStyle = 1
$inputTop.Interior.ColorIndex = 2
$inputTop.Font.ThemeColor = 1
$inputTop.Font.Name = "Arial"
$inputTop.VerticalAlignment = -4160
